$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "69.044.46"
Set-TextCell 2 5 "  +0.62%  "

Set-TextCell 3 4 "3.741.13"
Set-TextCell 3 5 "  +1.22%  "

Set-TextCell 4 5 "  +0.03%  "

Set-TextCell 5 4 "601.99"
Set-TextCell 5 5 "  +0.41%  "

Set-TextCell 6 4 "167.69"
Set-TextCell 6 5 "  +0.45%  "

Set-TextCell 7 4 "3.739.58"
Set-TextCell 7 5 "  +1.25%  "

Set-TextCell 8 5 "  -0.03%  "

Set-TextCell 9 5 "  +0.98%  "

Set-TextCell 10 5 "  +2.67%  "

Set-TextCell 11 4 "6.46"
Set-TextCell 11 5 "  +2.79%  "

Set-TextCell 12 5 "  +0.52%  "

Set-TextCell 13 4 "38.03"
Set-TextCell 13 5 "  +0.05%  "

Set-TextCell 14 4 "0.0000248"
Set-TextCell 14 5 "  +2.41%  "

Set-TextCell 15 4 "4.368.67"
Set-TextCell 15 5 "  +1.40%  "

Set-TextCell 16 4 "3.743.56"
Set-TextCell 16 5 "  +1.29%  "

Set-TextCell 17 4 "69.011.05"
Set-TextCell 17 5 "  +0.62%  "

Set-TextCell 18 4 "7.34"
Set-TextCell 18 5 "  +1.34%  "

Set-TextCell 19 2 "Chainlink"
Set-TextCell 19 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell 19 4 "17.28"
Set-TextCell 19 5 "  +1.12%  "

Set-TextCell 20 2 "TRON"
Set-TextCell 20 3 "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell 20 4 "0.113"
Set-TextCell 20 5 "  -1.26%  "

Set-TextCell 21 4 "10.90"
Set-TextCell 21 5 "  +19.57%  "

Set-TextCell 22 4 "493.13"
Set-TextCell 22 5 "  +0.34%  "

Set-TextCell 23 4 "0.725"
Set-TextCell 23 5 "  +0.66%  "

Set-TextCell 24 4 "0.0000153"
Set-TextCell 24 5 "  +8.55%  "

Set-TextCell 25 4 "84.79"
Set-TextCell 25 5 "  +0.46%  "

Set-TextCell 26 5 "  +0.74%  "

Set-TextCell 27 4 "12.36"
Set-TextCell 27 5 "  +1.50%  "

Set-TextCell 28 4 "10.09"
Set-TextCell 28 5 "  +0.37%  "

Set-TextCell 29 5 "  -0.07%  "

Set-TextCell 30 5 "  +2.35%  "

Set-TextCell 31 5 "  +4.77%  "

Set-TextCell 32 4 "8.05"
Set-TextCell 32 5 "  +3.36%  "

Set-TextCell 33 4 "31.51"
Set-TextCell 33 5 "  +0.34%  "

Set-TextCell 34 4 "3.886.78"
Set-TextCell 34 5 "  +1.56%  "

Set-TextCell 35 2 "Hedera"
Set-TextCell 35 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 35 4 "0.109"
Set-TextCell 35 5 "  +0.38%  "

Set-TextCell 36 2 "RenzoRestakedETH"
Set-TextCell 36 3 "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextCell 36 4 "3.677.81"
Set-TextCell 36 5 "  +1.17%  "

Set-TextCell 37 4 "0.999"
Set-TextCell 37 5 "  +0.04%  "

Set-TextCell 38 5 "  +1.68%  "

Set-TextCell 39 4 "5.93"
Set-TextCell 39 5 "  +3.59%  "

Set-TextCell 40 5 "  +1.92%  "

Set-TextCell 41 5 "  +1.04%  "

Set-TextCell 42 4 "2.99"
Set-TextCell 42 5 "  +6.37%  "

Set-TextCell 43 4 "431.68"
Set-TextCell 43 5 "  -0.05%  "

Set-TextCell 44 2 "Stacks"
Set-TextCell 44 3 "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell 44 4 "2.00"
Set-TextCell 44 5 "  +2.00%  "

Set-TextCell 45 2 "OKB"
Set-TextCell 45 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell 45 4 "48.64"
Set-TextCell 45 5 "  -0.81%  "

Set-TextCell 46 4 "8.48"
Set-TextCell 46 5 "  +1.35%  "

Set-TextCell 47 5 "  +0.01%  "

Set-TextCell 48 4 "40.38"
Set-TextCell 48 5 "  +0.49%  "

Set-TextCell 49 4 "141.02"
Set-TextCell 49 5 "  -0.32%  "

Set-TextCell 50 4 "2.782.00"
Set-TextCell 50 5 "  +1.89%  "

Set-TextCell 51 5 "  +0.93%  "
